$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of data to append to the timeline (rows 16-18), matching the
# formatting of the existing plain data cells (e.g. A15/C15 which use the
# default centered "General" style). Copy formatting from row 15's
# A/C columns (style index 1) across the new A16:D18 block so no new
# style entries are introduced.
$ws.Range("A15:A15").Copy() | Out-Null
$ws.Range("A16:D18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 16: Day 13 (entered completely, left to right)
$ws.Cells.Item(16, 1).Value = 13
$ws.Cells.Item(16, 2).Value = "14/3/2024"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = "refactored cache and Added User Basic Details Update"

# Day/Date/Hours for rows 17-18 filled in next
$ws.Cells.Item(17, 1).Value = 14
$ws.Cells.Item(17, 2).Value = "15/3/2024"
$ws.Cells.Item(17, 3).Value = 5

$ws.Cells.Item(18, 1).Value = 15
$ws.Cells.Item(18, 2).Value = "16/3/2024"
$ws.Cells.Item(18, 3).Value = 6

# Descriptions for rows 17-18 filled in last
$ws.Cells.Item(17, 4).Value = "Fight with cache second level"
$ws.Cells.Item(18, 4).Value = "Finally optimized the cache in a great way, refactored user specific detail"

# Recalculate so the SUM formula in D36 reflects the newly added hours.
$excel.Calculate()

# Update the selected cell to match the saved workbook view state.
$ws.Range("D22").Select()
